# Revisión de la Primera historia de Usuario
# Update the EXECUTED-sprint numbers for the first batch of user stories
# and move the active selection to the top of the backlog table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("EXECUTED") — revised sprint numbers after reviewing the
# first user story.
$ws.Range("F9").Value  = 4
$ws.Range("F10").Value = 6
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 6
$ws.Range("F13").Value = 5

# Leave the cursor parked on the backlog header area.
[void]$ws.Range("B4").Select()
